# relation and indexing problem v1.1
#
# The "BIGSERIES" data type used throughout the DB Architecture schema is
# renamed to "BIGINT 20". Every table's primary-key type cell (column A of
# each "Typy" row) currently reads "BIGSERIES"; replace all of them with
# "BIGINT 20" in one shot, then restore the sheet's scroll position / active
# selection to cell G8 (top of sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace every "BIGSERIES" cell value with "BIGINT 20" across the sheet.
$ws.Cells.Replace("BIGSERIES", "BIGINT 20", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# Scroll back to the top and select G8 (matches the saved view state).
$ws.Range("A1").Select()
$ws.Range("G8").Select()
